$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.765.36'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.091.74'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  -2.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.82%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.085.12'
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.157'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.35'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000227'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('D15').Value = '3.594.42'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '63.787.43'
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '3.093.76'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '489.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.82%  '
$ws.Range('E22').Value = '  -1.31%  '
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.89'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.30'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('E31').Value = '  -2.81%  '
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('E33').Value = '  -5.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '56.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '493.86'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.50%  '
$ws.Range('D38').Value = '3.315.94'
$ws.Range('E38').Value = '  +7.66%  '
$ws.Range('E39').Value = '  -4.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0801'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.260'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').Value = '0.0₃0540'
$ws.Range('E47').Value = '  +5.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('E51').Value = '  -14.98%  '
